# empty_question_label_patch.xlsx - "start" row split into two "text" rows
# plus the survey sheet becomes the active/selected tab (settings loses it).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# --- Insert a new row for the table so the "start" row becomes two rows ---
# Row 4 currently holds "end group"/"a"; pushing it down to row 5 opens up
# a fresh row 4 that inherits row 3's formatting (the row being split).
$ws.Rows.Item(4).Insert()

# Grow the "survey" table (ListObject) so it covers the new row too.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D5"))

# --- Re-purpose old row 3 ("start"/"start"/My image label/my_image1.jpg)
#     into the first "text" question, and populate the new row 4 with the
#     second "text" question. Write order matters for shared-string ids. ---
$ws.Range("A3").Value = "text"
$ws.Range("B3").Value = "item1"
$ws.Range("A4").Value = "text"
$ws.Range("B4").Value = "item2"
$ws.Range("C4").Value = "My image label 2"
$ws.Range("C3").Value = "My image label 1"
$ws.Range("D4").Value = "my_image2.jpg"
# D3 already holds "my_image1.jpg" - no change needed there.

# Match the "begin group"/"end group" row styling (font + text number format)
# on the newly-inserted row 4, and keep row 3 consistent across its columns.
$ws.Range("B4:C4").NumberFormat = "@"
$ws.Range("A3:C3").NumberFormat = "@"

# --- Make "survey" the selected/active sheet with D4 highlighted ---
$ws.Activate()
$ws.Range("D4").Select()

# --- "settings" should no longer be the tab-selected sheet ---
$ws2 = $wb.Worksheets.Item("settings")
$ws2Window = $excel.ActiveWindow
